# The shared string "(%, mIoU)" was removed from the workbook. It was
# only ever used in column R of rows 3, 7, and 18-37, so clearing those
# cells removes the last references to that shared string (which the
# writer then prunes from the shared-strings table on save, shifting all
# later shared-string indices down by one - this happens automatically,
# we just need to clear the cell contents).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R3").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("R18:R37").ClearContents()

# Reflect the final selected cell recorded in the saved workbook.
$ws.Range("R16").Select()
